$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("一些计算")
$ws.Range("A1").Value = "Test"
